$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trivial")

# Add a new row of data (row 6) following the same pattern as existing rows
$ws.Range("A6").Value = "preguntador 5"
$ws.Range("B6").Value = "pregunta 5"
$ws.Range("C6").Value = "respuesta 5"

# Column C uses a text-formatted style (style index 1) like the rows above it
$ws.Range("C6").NumberFormat = "@"

# Update the active selection to match the author's final cursor position
$ws.Range("E5").Select()
